# Adiciona as colunas de volume, area e comprimento por dm3 (por litro de solo)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Novos cabeçalhos (linha 1) com o mesmo estilo (negrito/centralizado) das demais colunas
$ws.Range("I1").Value = "volume_cm3_dm3"
$ws.Range("J1").Value = "area_cm2_dm3"
$ws.Range("K1").Value = "comprimento_cm_dm3"
$ws.Range("I1:K1").Font.Bold = $true
$ws.Range("I1:K1").HorizontalAlignment = -4108

# Calcula os valores por dm3 para cada linha de dados (2..última linha usada)
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $volume = $ws.Cells.Item($r, 5).Value()       # coluna E - volume
    $area = $ws.Cells.Item($r, 6).Value()          # coluna F - area
    $comprimento = $ws.Cells.Item($r, 8).Value()   # coluna H - comprimento

    $ws.Cells.Item($r, 9).Value = $volume / 475
    $ws.Cells.Item($r, 10).Value = $area / 47.5
    $ws.Cells.Item($r, 11).Value = $comprimento / 4.75
}
